$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 493 ("「手書きフォント」" entry) entirely; this shifts all
# subsequent rows up by one, matching the target workbook's row layout.
$ws.Rows("493:493").Delete()
